# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Updates column G ("K") values for rows 2-38 (except row 37, which is
# unchanged) on the active worksheet to reflect the recalculated
# strikeout (K) totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 9
    4  = 6
    5  = 4
    6  = 7
    7  = 4
    8  = 9
    9  = 12
    10 = 4
    11 = 8
    12 = 5
    13 = 11
    14 = 6
    15 = 3
    16 = 6
    17 = 5
    18 = 3
    19 = 7
    20 = 2
    21 = 7
    22 = 10
    23 = 3
    24 = 6
    25 = 11
    26 = 5
    27 = 7
    28 = 7
    29 = 4
    30 = 13
    31 = 6
    32 = 15
    33 = 8
    34 = 5
    35 = 5
    36 = 4
    38 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
